$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 140, shifting existing rows 140-146 down to 141-147
$ws.Rows.Item(140).Insert()

# Populate the new row 140 with the new data record
$ws.Cells.Item(140, 1).Value = 10
$ws.Cells.Item(140, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(140, 3).Value = "La Araucanía"
$ws.Cells.Item(140, 4).Value = 45021
$ws.Cells.Item(140, 5).Value = 9
$ws.Cells.Item(140, 6).Value = 100112035
$ws.Cells.Item(140, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 30
$ws.Cells.Item(140, 11).Value = 30000
$ws.Cells.Item(140, 12).Value = 30000
$ws.Cells.Item(140, 13).Value = 30000
$ws.Cells.Item(140, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(140, 15).Value = "Región Metropolitana"
$ws.Cells.Item(140, 16).Value = 2000
$ws.Cells.Item(140, 17).Value = 15
$ws.Cells.Item(140, 18).Value = "Hortaliza"
